$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '89.461.66'
$ws.Cells.Item(2, 5).Value = '  +10.65%  '

$ws.Cells.Item(3, 4).Value = '3.366.80'
$ws.Cells.Item(3, 5).Value = '  +7.53%  '

$ws.Cells.Item(4, 5).Value = '  +0.17%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '220.24'
$ws.Cells.Item(5, 5).Value = '  +6.39%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '650.31'
$ws.Cells.Item(6, 5).Value = '  +5.86%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.408'
$ws.Cells.Item(7, 5).Value = '  +47.22%  '

$ws.Cells.Item(8, 5).Value = '  +0.05%  '

$ws.Cells.Item(9, 5).Value = '  +7.14%  '

$ws.Cells.Item(10, 4).Value = '3.370.37'
$ws.Cells.Item(10, 5).Value = '  +7.66%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.611'
$ws.Cells.Item(11, 5).Value = '  +7.60%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.0000289'
$ws.Cells.Item(12, 5).Value = '  +16.43%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '36.69'
$ws.Cells.Item(13, 5).Value = '  +17.95%  '

$ws.Cells.Item(14, 5).Value = '  +2.47%  '

$ws.Cells.Item(15, 4).Value = '3.994.32'
$ws.Cells.Item(15, 5).Value = '  +7.80%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '5.58'
$ws.Cells.Item(16, 5).Value = '  +6.68%  '

$ws.Cells.Item(17, 4).Value = '89.351.44'
$ws.Cells.Item(17, 5).Value = '  +10.56%  '

$ws.Cells.Item(18, 4).Value = '3.367.20'
$ws.Cells.Item(18, 5).Value = '  +7.84%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '14.87'
$ws.Cells.Item(19, 5).Value = '  +8.33%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '3.17'
$ws.Cells.Item(20, 5).Value = '  +2.10%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '9.68'
$ws.Cells.Item(21, 5).Value = '  +9.40%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '458.69'
$ws.Cells.Item(22, 5).Value = '  +7.32%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '5.55'
$ws.Cells.Item(23, 5).Value = '  +10.35%  '

$ws.Cells.Item(24, 5).Value = '  +5.67%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '5.58'
$ws.Cells.Item(25, 5).Value = '  +9.64%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '12.87'
$ws.Cells.Item(26, 5).Value = '  +20.20%  '

$ws.Cells.Item(27, 4).Value = '3.524.32'
$ws.Cells.Item(27, 5).Value = '  +7.13%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '0.0000143'
$ws.Cells.Item(28, 5).Value = '  +20.24%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '79.13'
$ws.Cells.Item(29, 5).Value = '  +5.06%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.202'
$ws.Cells.Item(30, 5).Value = '  +45.77%  '

$ws.Cells.Item(31, 5).Value = '  -0.13%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '9.43'
$ws.Cells.Item(32, 5).Value = '  +6.94%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '598.18'
$ws.Cells.Item(33, 5).Value = '  +8.60%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.58'
$ws.Cells.Item(34, 5).Value = '  +8.86%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.00'
$ws.Cells.Item(35, 5).Value = '  -0.01%  '

$ws.Cells.Item(36, 5).Value = '  +8.00%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '7.42'
$ws.Cells.Item(37, 5).Value = '  +25.07%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.145'
$ws.Cells.Item(38, 5).Value = '  -3.39%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '23.51'
$ws.Cells.Item(39, 5).Value = '  +5.09%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.430'
$ws.Cells.Item(40, 5).Value = '  +7.31%  '

$ws.Cells.Item(41, 5).Value = '  +9.04%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '3.22'
$ws.Cells.Item(42, 5).Value = '  +8.22%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '21.87'
$ws.Cells.Item(43, 5).Value = '  +5.61%  '

$ws.Cells.Item(44, 5).Value = '  +0.07%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '1.48'
$ws.Cells.Item(45, 5).Value = '  +13.95%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '158.00'
$ws.Cells.Item(46, 5).Value = '  -0.53%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '190.70'
$ws.Cells.Item(48, 5).Value = '  +2.69%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '46.39'
$ws.Cells.Item(49, 5).Value = '  +4.14%  '

$ws.Cells.Item(50, 5).Value = '  +9.37%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.672'
$ws.Cells.Item(51, 5).Value = '  +8.84%  '
